$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update existing rows 2-5 with new schedule values
$data = @(
    @(1, 5, 3, 3, 7, -2, 4, 45, 5),
    @(2, 6, 2, 5, 7, -1, 5, 56, 5),
    @(3, 6, 4, 1, 5, -5, 1, 12, 5),
    @(4, 5, 1, 2, 4, -3, 3, 34, 5),
    @(5, 8, 4, 4, 6, -4, 2, 23, 5)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    for ($c = 0; $c -lt $vals.Length; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
    }
    $ws.Cells.Item($row, 10).Value = "train_dim2_1"
}

# Update the selected cell to I1
$ws.Range("I1").Select()
